$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: a numeric value (stored/formatted as text) in column B
# and a text value in column C, used when searching a column for a value.
$ws.Range("B8").Value = 12
$ws.Range("B8").NumberFormat = "@"
$ws.Range("C8").Value = "wer"

[void]$ws.Range("B8").Select()
